$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 252; this shifts existing rows 252..283 down to 253..284
$ws.Rows(252).Insert()

# Populate the newly inserted row 252 with the new weekly price record.
# Columns A,B,C,E,F,G,H,I,J,K,R are constant across this product's rows.
$ws.Cells.Item(252, 1).Value = 10
$ws.Cells.Item(252, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(252, 3).Value = "La Araucanía"
$ws.Cells.Item(252, 4).Value = 45077
$ws.Cells.Item(252, 5).Value = 9
$ws.Cells.Item(252, 6).Value = "Fruta"
$ws.Cells.Item(252, 7).Value = 100104
$ws.Cells.Item(252, 8).Value = "Frutos de pepita"
$ws.Cells.Item(252, 9).Value = 100104003
$ws.Cells.Item(252, 10).Value = "Membrillo"
$ws.Cells.Item(252, 11).Value = "Champion"
$ws.Cells.Item(252, 12).Value = "Primera"
$ws.Cells.Item(252, 13).Value = 80
$ws.Cells.Item(252, 14).Value = 14000
$ws.Cells.Item(252, 15).Value = 14000
$ws.Cells.Item(252, 16).Value = 14000
$ws.Cells.Item(252, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(252, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(252, 19).Value = 778
$ws.Cells.Item(252, 20).Value = 18
